$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.739.56'
$ws.Range('E2').Value = '  -1.74%  '
$ws.Range('D3').Value = '1.759.32'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.21%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4429'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.90%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3743'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.48'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07674'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.76'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.196'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.433'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '1.757.55'
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001078'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.24'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +9.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06229'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.61%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.180'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5327'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.19%  '
$ws.Range('D24').Value = '27.761.95'
$ws.Range('E24').Value = '  -1.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.65'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.34%  '
$ws.Range('E26').Value = '  -4.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.78'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '153.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.359'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').Value = '1.956.20'
$ws.Range('E30').Value = '  -2.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '128.48'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.218'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.81%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.767'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09315'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.651'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -9.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.71'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2191'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02323'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06150'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6496'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.092'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.201'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.419'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.29%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6023'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.766'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '126.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.001'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.140'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.74%  '
